$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2433.923
$ws.Range("I15").Value = 2433.923
$ws.Range("K15").Value = 7301.768999999999
$ws.Range("M15").Value = -7132.768999999999
$ws.Range("H33").Value = 5439.05
$ws.Range("I33").Value = 6134.353
$ws.Range("K33").Value = 6134.353
$ws.Range("M33").Value = -5905.353
$ws.Range("H62").Value = 12655.389
$ws.Range("J62").Value = 9658.923000000001
$ws.Range("L62").Value = 9658.923000000001
$ws.Range("N62").Value = -10906.923
$ws.Range("H65").Value = 12655.389
$ws.Range("J65").Value = 9658.923000000001
$ws.Range("L65").Value = 48294.61500000001
$ws.Range("N65").Value = -54534.61500000001
$ws.Range("H111").Value = 1407.8889
$ws.Range("I111").Value = 1407.8889
$ws.Range("K111").Value = 4223.6667
$ws.Range("M111").Value = -1156.6667
$ws.Range("H116").Value = 4982.864
$ws.Range("I116").Value = 4444.7144
$ws.Range("K116").Value = 4444.7144
$ws.Range("M116").Value = -1002.7144
$ws.Range("H132").Value = 9166.261
$ws.Range("I132").Value = 8062.488
$ws.Range("K132").Value = 24187.464
$ws.Range("M132").Value = -21657.464
$ws.Range("H137").Value = 1712.7
$ws.Range("I137").Value = 890
$ws.Range("J137").Value = 3632.3333
$ws.Range("K137").Value = 2670
$ws.Range("L137").Value = 10896.9999
$ws.Range("M137").Value = -120
$ws.Range("N137").Value = -15996.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 522.4
$ws.Range("I2").Value = 586.125
$ws.Range("J2").Value = 267.5
$ws.Range("K2").Value = 586.125
$ws.Range("L2").Value = 267.5
$ws.Range("M2").Value = -473.125
$ws.Range("N2").Value = -493.5
$ws.Range("H32").Value = 24005.408
$ws.Range("I32").Value = 24054.844
$ws.Range("K32").Value = 24054.844
$ws.Range("M32").Value = -23767.844
$ws.Range("H116").Value = 522.4
$ws.Range("I116").Value = 586.125
$ws.Range("J116").Value = 267.5
$ws.Range("K116").Value = 586.125
$ws.Range("L116").Value = 267.5
$ws.Range("M116").Value = 1707.875
$ws.Range("N116").Value = -4855.5
$ws.Range("H122").Value = 1801.0588
$ws.Range("I122").Value = 1561.2667
$ws.Range("K122").Value = 4683.800099999999
$ws.Range("M122").Value = -2233.800099999999
$ws.Range("H132").Value = 23586.83
$ws.Range("I132").Value = 26013.977
$ws.Range("J132").Value = 3198.8
$ws.Range("K132").Value = 78041.931
$ws.Range("L132").Value = 9596.400000000001
$ws.Range("M132").Value = -75511.931
$ws.Range("N132").Value = -14656.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 522.4
$ws.Range("I3").Value = 586.125
$ws.Range("J3").Value = 267.5
$ws.Range("K3").Value = 586.125
$ws.Range("L3").Value = 267.5
$ws.Range("M3").Value = -472.125
$ws.Range("N3").Value = -495.5
$ws.Range("H105").Value = 4773
$ws.Range("I105").Value = 4355.1177
$ws.Range("J105").Value = 4970.3335
$ws.Range("K105").Value = 4355.1177
$ws.Range("L105").Value = 4970.3335
$ws.Range("M105").Value = -2608.1177
$ws.Range("N105").Value = -8464.333500000001
$ws.Range("H134").Value = 2483.0312
$ws.Range("I134").Value = 2160.9656
$ws.Range("K134").Value = 6482.8968
$ws.Range("M134").Value = -3947.8968

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 83.2
$ws.Range("I7").Value = 64.666664
$ws.Range("J7").Value = 111
$ws.Range("K7").Value = 64.666664
$ws.Range("L7").Value = 111
$ws.Range("M7").Value = 48.333336
$ws.Range("N7").Value = -337
$ws.Range("H31").Value = 2540.4062
$ws.Range("I31").Value = 2276.087
$ws.Range("J31").Value = 3215.889
$ws.Range("K31").Value = 2276.087
$ws.Range("L31").Value = 3215.889
$ws.Range("M31").Value = -1981.087
$ws.Range("N31").Value = -3805.889
$ws.Range("H34").Value = 2540.4062
$ws.Range("I34").Value = 2276.087
$ws.Range("J34").Value = 3215.889
$ws.Range("K34").Value = 2276.087
$ws.Range("L34").Value = 3215.889
$ws.Range("M34").Value = -2074.087
$ws.Range("N34").Value = -3619.889
$ws.Range("H105").Value = 1299
$ws.Range("I105").Value = 1299
$ws.Range("K105").Value = 1299
$ws.Range("M105").Value = 448
$ws.Range("H107").Value = 473.5263
$ws.Range("I107").Value = 455.70587
$ws.Range("K107").Value = 455.70587
$ws.Range("M107").Value = 1464.29413

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1194804.6
$ws.Range("I4").Value = 1736562.6
$ws.Range("J4").Value = 2937
$ws.Range("K4").Value = 5209687.800000001
$ws.Range("L4").Value = 8811
$ws.Range("M4").Value = -5209575.800000001
$ws.Range("N4").Value = -9035
$ws.Range("H62").Value = 8175.6665
$ws.Range("J62").Value = 8175.6665
$ws.Range("L62").Value = 24526.9995
$ws.Range("N62").Value = -25898.9995
$ws.Range("H65").Value = 8175.6665
$ws.Range("J65").Value = 8175.6665
$ws.Range("L65").Value = 73580.9985
$ws.Range("N65").Value = -80444.9985
$ws.Range("H75").Value = 777
$ws.Range("I75").Value = 777
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 2331
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -1333
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 777
$ws.Range("I78").Value = 777
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 6993
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -2001
$ws.Range("N78").ClearContents()
$ws.Range("H113").Value = 917.86957
$ws.Range("I113").Value = 952.6667
$ws.Range("K113").Value = 2858.0001
$ws.Range("M113").Value = -688.0001000000002
$ws.Range("H114").Value = 1194.3334
$ws.Range("I114").Value = 791.5
$ws.Range("J114").Value = 2000
$ws.Range("K114").Value = 2374.5
$ws.Range("L114").Value = 6000
$ws.Range("M114").Value = 879.5
$ws.Range("N114").Value = -12508
$ws.Range("H137").Value = 4166
$ws.Range("I137").Value = 1499.5
$ws.Range("K137").Value = 4498.5
$ws.Range("M137").Value = 601.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14458.333
$ws.Range("I70").Value = 14687.75
$ws.Range("K70").Value = 14687.75
$ws.Range("M70").Value = -14417.75
$ws.Range("H73").Value = 14458.333
$ws.Range("I73").Value = 14687.75
$ws.Range("K73").Value = 14687.75
$ws.Range("M73").Value = -13751.75
$ws.Range("H102").Value = 2007.4615
$ws.Range("I102").Value = 779.55554
$ws.Range("K102").Value = 779.55554
$ws.Range("M102").Value = 842.44446
$ws.Range("H113").Value = 422000.2
$ws.Range("I113").Value = 368333.66
$ws.Range("K113").Value = 368333.66
$ws.Range("M113").Value = -366163.66
$ws.Range("H122").Value = 2787.5
$ws.Range("I122").Value = 2391
$ws.Range("K122").Value = 7173
$ws.Range("M122").Value = -4723

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 9496.529
$ws.Range("I46").Value = 19219.6
$ws.Range("K46").Value = 19219.6
$ws.Range("M46").Value = -19031.6
$ws.Range("H63").Value = 51999.5
$ws.Range("J63").Value = 63999
$ws.Range("L63").Value = 63999
$ws.Range("N63").Value = -65497
$ws.Range("H66").Value = 51999.5
$ws.Range("J66").Value = 63999
$ws.Range("L66").Value = 191997
$ws.Range("N66").Value = -199485

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 742.1579
$ws.Range("I107").Value = 725.1177
$ws.Range("J107").Value = 887
$ws.Range("K107").Value = 2175.3531
$ws.Range("L107").Value = 2661
$ws.Range("M107").Value = -255.3531000000003
$ws.Range("N107").Value = -6501
$ws.Range("H113").Value = 1150.2174
$ws.Range("I113").Value = 1127.6
$ws.Range("K113").Value = 3382.8
$ws.Range("M113").Value = -1212.8
$ws.Range("H124").Value = 51749.75
$ws.Range("J124").Value = 51749.75
$ws.Range("L124").Value = 51749.75
$ws.Range("N124").Value = -61569.75
$ws.Range("H132").Value = 60308.6
$ws.Range("I132").Value = 60308.6
$ws.Range("K132").Value = 180925.8
$ws.Range("M132").Value = -178395.8
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
